# Issue #13: Permitir que en los ficheros de metadatos dos columnas se puedan
# relacionar para crear SKOS jerarquicos.
#
# A new row is inserted right below the header row. It holds the
# "slug"/identifier form of each column header (used later to relate two
# columns hierarchically), pushing the former rows 2-4 (URIs, medida/dim,
# datatypes) down to rows 3-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2, shifting existing rows 2..4 down to 3..4.
$ws.Rows("2:2").Insert()

# Fill the new row with the slugified header identifiers.
$ws.Range("A2").Value = "numero-hogares"
$ws.Range("B2").Value = "municipio-2-residencia-nombre"
$ws.Range("C2").Value = "aragon"
$ws.Range("D2").Value = "municipio-2-residencia-codigo"
